$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1751
$ws.Range("I12").Value = 1500
$ws.Range("J12").Value = 2002
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 2002
$ws.Range("M12").Value = -1330
$ws.Range("N12").Value = -2342

# Row 74
$ws.Range("H74").Value = 3549.5
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 3599
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 3599
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -5471

# Row 77
$ws.Range("H77").Value = 3549.5
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 3599
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 17995
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -27355

# Row 88
$ws.Range("H88").Value = 426.41666
$ws.Range("I88").Value = 349.75
$ws.Range("J88").Value = 464.75
$ws.Range("K88").Value = 349.75
$ws.Range("L88").Value = 464.75
$ws.Range("M88").Value = 56.25
$ws.Range("N88").Value = -1276.75

# Row 91
$ws.Range("H91").Value = 426.41666
$ws.Range("I91").Value = 349.75
$ws.Range("J91").Value = 464.75
$ws.Range("K91").Value = 349.75
$ws.Range("L91").Value = 464.75
$ws.Range("M91").Value = 1054.25
$ws.Range("N91").Value = -3272.75

# Row 135
$ws.Range("H135").Value = 940.8461
$ws.Range("I135").Value = 757.36365
$ws.Range("J135").Value = 1950
$ws.Range("K135").Value = 6816.27285
$ws.Range("L135").Value = 17550
$ws.Range("M135").Value = -4281.27285

# Row 137
$ws.Range("H137").Value = 4817.636
$ws.Range("I137").Value = 1998.7778
$ws.Range("J137").Value = 6769.154
$ws.Range("K137").Value = 5996.3334
$ws.Range("L137").Value = 20307.462
$ws.Range("M137").Value = -3446.3334


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 405
$ws.Range("I3").Value = 405
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 405
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -290
$ws.Range("N3").ClearContents()

# Row 32
$ws.Range("H32").Value = 2947.3713
$ws.Range("I32").Value = 2823
$ws.Range("J32").Value = 4999.5
$ws.Range("K32").Value = 2823
$ws.Range("L32").Value = 4999.5
$ws.Range("M32").Value = -2536

# Row 74
$ws.Range("H74").Value = 1881.6666
$ws.Range("I74").Value = 1475.75
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 1475.75
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -601.75

# Row 77
$ws.Range("H77").Value = 1881.6666
$ws.Range("I77").Value = 1475.75
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 7378.75
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -3010.75

# Row 132
$ws.Range("H132").Value = 2791.0688
$ws.Range("I132").Value = 2365.3684
$ws.Range("J132").Value = 3599.9
$ws.Range("K132").Value = 7096.1052
$ws.Range("L132").Value = 10799.7
$ws.Range("M132").Value = -4566.1052


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 590
$ws.Range("I12").Value = 621.6667
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 621.6667
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = -453.6667
$ws.Range("N12").Value = -736

# Row 20
$ws.Range("H20").Value = 3493
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3493
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3493
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -3987

# Row 86
$ws.Range("H86").Value = 6970
$ws.Range("I86").Value = 4697.5
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 4697.5
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -3574.5

# Row 89
$ws.Range("H89").Value = 6970
$ws.Range("I89").Value = 4697.5
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 23487.5
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -17871.5

# Row 94
$ws.Range("H94").Value = 1296.7142
$ws.Range("I94").Value = 1351.3334
$ws.Range("J94").Value = 969
$ws.Range("K94").Value = 1351.3334
$ws.Range("L94").Value = 969
$ws.Range("M94").Value = -900.3334


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 6301.125
$ws.Range("I94").Value = 3276.4
$ws.Range("J94").Value = 11342.333
$ws.Range("K94").Value = 3276.4
$ws.Range("L94").Value = 11342.333
$ws.Range("M94").Value = -2825.4

# Row 132
$ws.Range("H132").Value = 3517
$ws.Range("I132").Value = 2896.25
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 8688.75
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -6158.75


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1381.125
$ws.Range("I5").Value = 1175
$ws.Range("J5").Value = 1587.25
$ws.Range("K5").Value = 3525
$ws.Range("L5").Value = 4761.75
$ws.Range("M5").Value = -3413
$ws.Range("N5").Value = -4985.75

# Row 12
$ws.Range("H12").Value = 108.52941
$ws.Range("I12").Value = 12
$ws.Range("J12").Value = 148.75
$ws.Range("K12").Value = 36
$ws.Range("L12").Value = 446.25
$ws.Range("M12").Value = 137

# Row 13
$ws.Range("H13").Value = 308.33334
$ws.Range("I13").Value = 308.33334
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 925.0000200000001
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -757.0000200000001
$ws.Range("N13").ClearContents()

# Row 17
$ws.Range("H17").Value = 2694.2
$ws.Range("I17").Value = 49.8
$ws.Range("J17").Value = 3575.6667
$ws.Range("K17").Value = 149.4
$ws.Range("L17").Value = 10727.0001
$ws.Range("M17").Value = 19.60000000000002
$ws.Range("N17").Value = -11065.0001

# Row 135
$ws.Range("H135").Value = 1381.125
$ws.Range("I135").Value = 1175
$ws.Range("J135").Value = 1587.25
$ws.Range("K135").Value = 10575
$ws.Range("L135").Value = 14285.25
$ws.Range("M135").Value = -8040
$ws.Range("N135").Value = -19355.25


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 95
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 95
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 17

# Row 70
$ws.Range("H70").Value = 4999
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4999
$ws.Range("N70").Value = -5539

# Row 73
$ws.Range("H73").Value = 4999
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4999
$ws.Range("N73").Value = -6871

# Row 97
$ws.Range("H97").Value = 677.6667
$ws.Range("I97").Value = 677.6667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 677.6667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -181.6667

# Row 113
$ws.Range("H113").Value = 7791.3
$ws.Range("I113").Value = 6739
$ws.Range("J113").Value = 8843.6
$ws.Range("K113").Value = 6739
$ws.Range("L113").Value = 8843.6
$ws.Range("M113").Value = -4569

# Row 134
$ws.Range("H134").Value = 83979.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 83979.336
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 251938.008
$ws.Range("N134").Value = -257008.008


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2297.7778
$ws.Range("I132").Value = 2147.5
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6442.5
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3912.5
$ws.Range("N132").Value = -15560


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3821.4
$ws.Range("I81").Value = 4661.3335
$ws.Range("J81").Value = 2561.5
$ws.Range("K81").Value = 9322.666999999999
$ws.Range("L81").Value = 5123
$ws.Range("M81").Value = -8261.666999999999

# Row 84
$ws.Range("H84").Value = 3821.4
$ws.Range("I84").Value = 4661.3335
$ws.Range("J84").Value = 2561.5
$ws.Range("K84").Value = 46613.335
$ws.Range("L84").Value = 25615
$ws.Range("M84").Value = -41309.335

# Row 135
$ws.Range("H135").Value = 40000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

# Row 136
$ws.Range("H136").Value = 2697.68
$ws.Range("I136").Value = 1649.5294
$ws.Range("J136").Value = 4925
$ws.Range("K136").Value = 4948.5882
$ws.Range("L136").Value = 14775
$ws.Range("M136").Value = -2398.5882

